$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade record appended as row 16 (mirrors the existing rows' layout:
# Principle, Start Principle, BuyPrice, SellPrice, IsShortSell, Price Change %, Date, Profitable)
$ws.Range("A16").Value = 8561.7099999999991
$ws.Range("B16").Value = 8979.24
$ws.Range("C16").Value = 17.2
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = 4.6500000000000004
$ws.Range("G16").Value = 42626.545578703706
$ws.Range("H16").Value = $false

# Match the date formatting used by the rest of column G (style index 1)
# by copying the format from the row above instead of inventing a new
# number format entry.
$ws.Range("G15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
